# Scheduled-runner price refresh: overwrite computed market-price / profit
# columns (H..N) on affected leve rows across the ALC/ARM/BSM/CRP/CUL/LTW/WVR
# sheets with freshly pulled averages. Pure data overwrite - no formulas.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 6083.8335
$ws.Range("I64").Value = 6083.8335
$ws.Range("K64").Value = 6083.8335
$ws.Range("M64").Value = -5835.8335
$ws.Range("H67").Value = 6083.8335
$ws.Range("I67").Value = 6083.8335
$ws.Range("K67").Value = 6083.8335
$ws.Range("M67").Value = -5225.8335
$ws.Range("H74").Value = 11357.143
$ws.Range("I74").Value = 11357.143
$ws.Range("K74").Value = 11357.143
$ws.Range("M74").Value = -10421.143
$ws.Range("H77").Value = 11357.143
$ws.Range("I77").Value = 11357.143
$ws.Range("K77").Value = 56785.715
$ws.Range("M77").Value = -52105.715
$ws.Range("H86").Value = 4604
$ws.Range("I86").Value = 4808.8
$ws.Range("J86").Value = 4399.2
$ws.Range("K86").Value = 4808.8
$ws.Range("L86").Value = 4399.2
$ws.Range("M86").Value = -3685.8
$ws.Range("N86").Value = -6645.2
$ws.Range("H89").Value = 4604
$ws.Range("I89").Value = 4808.8
$ws.Range("J89").Value = 4399.2
$ws.Range("K89").Value = 24044
$ws.Range("L89").Value = 21996
$ws.Range("M89").Value = -18428
$ws.Range("N89").Value = -33228
$ws.Range("H132").Value = 2846.5945
$ws.Range("I132").Value = 2870.111
$ws.Range("K132").Value = 8610.332999999999
$ws.Range("M132").Value = -6080.332999999999
$ws.Range("H138").Value = 2606.5322
$ws.Range("J138").Value = 3424.139
$ws.Range("L138").Value = 10272.417
$ws.Range("N138").Value = -20552.417

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2732.257
$ws.Range("I32").Value = 1151.9181
$ws.Range("K32").Value = 1151.9181
$ws.Range("M32").Value = -864.9181000000001
$ws.Range("H63").Value = 2523.3333
$ws.Range("I63").Value = 2552.75
$ws.Range("J63").Value = 2499.8
$ws.Range("K63").Value = 2552.75
$ws.Range("L63").Value = 2499.8
$ws.Range("M63").Value = -1866.75
$ws.Range("N63").Value = -3871.8
$ws.Range("H66").Value = 2523.3333
$ws.Range("I66").Value = 2552.75
$ws.Range("J66").Value = 2499.8
$ws.Range("K66").Value = 12763.75
$ws.Range("L66").Value = 12499
$ws.Range("M66").Value = -9331.75
$ws.Range("N66").Value = -19363
$ws.Range("H97").Value = 494.86667
$ws.Range("I97").Value = 510.16666
$ws.Range("K97").Value = 510.16666
$ws.Range("M97").Value = -14.16665999999998
$ws.Range("H102").Value = 964
$ws.Range("I102").Value = 964
$ws.Range("K102").Value = 964
$ws.Range("M102").Value = 658
$ws.Range("H106").Value = 12499.5
$ws.Range("J106").Value = 12499.5
$ws.Range("L106").Value = 12499.5
$ws.Range("N106").Value = -15023.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 28337972
$ws.Range("J86").Value = 6053.5713
$ws.Range("L86").Value = 6053.5713
$ws.Range("N86").Value = -8299.5713
$ws.Range("H89").Value = 28337972
$ws.Range("J89").Value = 6053.5713
$ws.Range("L89").Value = 30267.8565
$ws.Range("N89").Value = -41499.85649999999
$ws.Range("H94").Value = 2465.0417
$ws.Range("I94").Value = 1266.3846
$ws.Range("J94").Value = 3881.6365
$ws.Range("K94").Value = 1266.3846
$ws.Range("L94").Value = 3881.6365
$ws.Range("M94").Value = -815.3846000000001
$ws.Range("N94").Value = -4783.636500000001
$ws.Range("H107").Value = 2325.6
$ws.Range("I107").Value = 1949.6666
$ws.Range("K107").Value = 1949.6666
$ws.Range("M107").Value = -29.66660000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2426.8333
$ws.Range("I62").Value = 2240.25
$ws.Range("K62").Value = 2240.25
$ws.Range("M62").Value = -1616.25
$ws.Range("H65").Value = 2426.8333
$ws.Range("I65").Value = 2240.25
$ws.Range("K65").Value = 11201.25
$ws.Range("M65").Value = -8081.25
$ws.Range("H105").Value = 1305.5
$ws.Range("I105").Value = 1363.25
$ws.Range("J105").Value = 1190
$ws.Range("K105").Value = 1363.25
$ws.Range("L105").Value = 1190
$ws.Range("M105").Value = 383.75
$ws.Range("N105").Value = -4684
$ws.Range("H134").Value = 811.2368
$ws.Range("I134").Value = 806.13513
$ws.Range("K134").Value = 2418.40539
$ws.Range("M134").Value = 116.5946100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 107.210526
$ws.Range("I38").Value = 109
$ws.Range("J38").Value = 104.75
$ws.Range("K38").Value = 327
$ws.Range("L38").Value = 314.25
$ws.Range("M38").Value = 20
$ws.Range("N38").Value = -1008.25
$ws.Range("H113").Value = 2300
$ws.Range("I113").Value = 450.5
$ws.Range("K113").Value = 1351.5
$ws.Range("M113").Value = 818.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 29017
$ws.Range("I7").Value = 29017
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 29017
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -28905
$ws.Range("N7").ClearContents()
$ws.Range("H122").Value = 9220.714
$ws.Range("I122").Value = 9717.947
$ws.Range("K122").Value = 29153.841
$ws.Range("M122").Value = -26703.841
$ws.Range("H126").Value = 29017
$ws.Range("I126").Value = 29017
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 87051
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -84581
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 200007
$ws.Range("J18").Value = 200007
$ws.Range("L18").Value = 200007
$ws.Range("N18").Value = -200353
$ws.Range("H26").Value = 20012
$ws.Range("I26").Value = 20012
$ws.Range("K26").Value = 20012
$ws.Range("M26").Value = -19719
$ws.Range("H96").Value = 812.25
$ws.Range("I96").Value = 795.0909
$ws.Range("J96").Value = 1001
$ws.Range("K96").Value = 795.0909
$ws.Range("L96").Value = 1001
$ws.Range("M96").Value = 577.9091
$ws.Range("N96").Value = -3747
$ws.Range("H100").Value = 1531.7368
$ws.Range("I100").Value = 1191.4
$ws.Range("J100").Value = 1909.8889
$ws.Range("K100").Value = 2382.8
$ws.Range("L100").Value = 3819.7778
$ws.Range("M100").Value = -1841.8
$ws.Range("N100").Value = -4901.7778
$ws.Range("H122").Value = 2115.5557
$ws.Range("I122").Value = 2019.3043
$ws.Range("K122").Value = 6057.9129
$ws.Range("M122").Value = -3607.9129
